$wb = $excel.ActiveWorkbook

# --- pitfall_data: rename ant-species headers to the new ".spN" naming ---
# (column positions are unchanged; only the species labels are updated)
$pitfall = $wb.Worksheets.Item("pitfall_data")
$pitfall.Range("F1").Value = "monomorium.sp1"
$pitfall.Range("G1").Value = "monomorium.sp2"
$pitfall.Range("K1").Value = "lepisiota.sp1"
$pitfall.Range("L1").Value = "lepisiota.sp2"
$pitfall.Range("H1").Value = "pheidole.sp1"
$pitfall.Range("I1").Value = "pheidole.sp2"
$pitfall.Range("J1").Value = "pheidole.sp3"

# --- metadata: replace the collector "Slie" with "Lindiwe Khoza" ---
$metadata = $wb.Worksheets.Item("metadata")
$metadata.Range("A10").VerticalAlignment = -4108   # xlCenter, creates the centered style
$metadata.Range("A10").Value = "Lindiwe Khoza"

# --- restore the view/selection state ---
$metadata.Activate()
$metadata.Range("J11").Select()

$pitfall.Activate()
$pitfall.Range("B1").Select()
